{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = new Set([\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n]);\n\n// Find the \"SERAFINI...\" paragraph, then delete the empty paragraph right\n// after it plus the two following text paragraphs (\"Ver no Jupiter...\" and\n// \"\u00a9 2020...\").\nlet serafiniIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"SERAFINI, Maria Jos\u00e9\") !== -1) {\n    serafiniIndex = i;\n    break;\n  }\n}\n\nif (serafiniIndex !== -1) {\n  const toDelete = [];\n  // the empty paragraph immediately following SERAFINI's paragraph\n  if (paragraphs.items[serafiniIndex + 1] && paragraphs.items[serafiniIndex + 1].text === \"\") {\n    toDelete.push(paragraphs.items[serafiniIndex + 1]);\n  }\n  for (let i = serafiniIndex + 2; i < paragraphs.items.length; i++) {\n    const text = paragraphs.items[i].text;\n    if (targets.has(text)) {\n      toDelete.push(paragraphs.items[i]);\n    }\n  }\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$copyrightText = \"$([char]0x00A9) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n$targets = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    $copyrightText\n)\n\n# Locate the \"SERAFINI...\" paragraph.\n$serafiniIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]0x0D, [char]0x07)\n    if ($text -like \"SERAFINI, Maria Jos*\") {\n        $serafiniIndex = $i\n        break\n    }\n}\n\nif ($serafiniIndex -ne -1) {\n    # Collect the indices to remove: the blank paragraph right after\n    # SERAFINI's paragraph, plus the \"Ver no Jupiter...\" and \"(c) 2020...\"\n    # paragraphs that follow it.\n    $indicesToDelete = New-Object System.Collections.ArrayList\n\n    $nextIndex = $serafiniIndex + 1\n    $nextText = $d.Paragraphs.Item($nextIndex).Range.Text.TrimEnd([char]0x0D, [char]0x07)\n    if ($nextText -eq \"\") {\n        [void]$indicesToDelete.Add($nextIndex)\n    }\n\n    for ($i = $serafiniIndex + 2; $i -le $d.Paragraphs.Count; $i++) {\n        $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]0x0D, [char]0x07)\n        if ($targets -contains $text) {\n            [void]$indicesToDelete.Add($i)\n        }\n    }\n\n    # Delete from the bottom up so earlier deletions don't shift the\n    # positions of paragraphs still queued for removal. Each paragraph is\n    # removed together with its trailing paragraph mark by spanning from\n    # its own start to the start of the following paragraph.\n    $sorted = $indicesToDelete | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $p = $d.Paragraphs.Item($idx)\n        $start = $p.Range.Start\n        if ($idx -lt $d.Paragraphs.Count) {\n            $end = $d.Paragraphs.Item($idx + 1).Range.Start\n        } else {\n            $end = $d.Content.End\n        }\n        $d.Range($start, $end).Delete()\n    }\n}\n"}
